$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. The "comments" column (M) only ever held a placeholder "NA" for
#    rows where no real comment existed. Those placeholders are being
#    dropped, which also removes the now-unused "NA" shared string.
# ------------------------------------------------------------------
$naRows = @(5, 6, 7, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23, 24, 25, 26, 27, 28)
foreach ($r in $naRows) {
    $ws.Cells.Item($r, 13).ClearContents()
}

# ------------------------------------------------------------------
# 2. Row 15 no longer carries a "group" label.
# ------------------------------------------------------------------
$ws.Range("B15").ClearContents()

# ------------------------------------------------------------------
# 3. New study "Morettini, 2019" is appended as row 29.
# ------------------------------------------------------------------
$ws.Range("A29").Value = "Morettini, 2019"
$ws.Range("C29").Value = 99
$ws.Range("D29").Value = 99
$ws.Range("E29").Value = 1346
$ws.Range("F29").Value = 0.005
$ws.Range("G29").Value = -0.49
$ws.Range("H29").Value = 0.5
$ws.Range("I29").Value = "Yes"
$ws.Range("J29").Value = "Eso"
$ws.Range("K29").Value = "Yes"
$ws.Range("L29").Value = "Elective Surgery"
$ws.Range("N29").Value = "Intraoperative"
$ws.Range("O29").Value = "low"
$ws.Range("P29").Value = "low"
$ws.Range("Q29").Value = "low"
$ws.Range("R29").Value = "low"
$ws.Range("S29").Value = "No"

# Wrap-text formatting for the numeric bias/CI columns, matching the
# rest of the table.
$ws.Range("F29:H29").WrapText = $true

# O29 picks up the same "no borders, alternate font" look already used
# by several RoB cells elsewhere in the column.
$ws.Range("O5").Copy()
$ws.Range("O29").PasteSpecial(-4122)
$ws.Range("O29").Value = "low"

$ws.Range("A1").Select()

# Move the active cell/selection to the newly-entered row, mirroring
# where the author's cursor ended up after adding the study.
$ws.Range("A29").Select()
